$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 22.77954133333333
$ws.Cells.Item(2, 8).Value = 68.338624
$ws.Cells.Item(2, 9).Value = 0.8649343844704168
$ws.Cells.Item(2, 10).Value = 0.8649343844704167
$ws.Cells.Item(2, 13).Value = 1.174933333333333
$ws.Cells.Item(2, 14).Value = 3.5248
$ws.Cells.Item(2, 15).Value = 0.01171850713626266
$ws.Cells.Item(2, 16).Value = 0.01171850713626266
$ws.Cells.Item(2, 17).Value = 26.76444243057778
$ws.Cells.Item(2, 18).Value = 240.8799818752
$ws.Cells.Item(2, 19).Value = 0.01013573975681553
$ws.Cells.Item(2, 20).Value = 0.01013573975681552
$ws.Cells.Item(3, 7).Value = 22.77954133333333
$ws.Cells.Item(3, 8).Value = 68.338624
$ws.Cells.Item(3, 9).Value = 0.8649343844704168
$ws.Cells.Item(3, 10).Value = 0.8649343844704167
$ws.Cells.Item(3, 15).Value = 0.2743256641287217
$ws.Cells.Item(3, 16).Value = 0.2743256641287218
$ws.Cells.Item(3, 17).Value = 626.5451187108124
$ws.Cells.Item(3, 18).Value = 5638.906068397312
$ws.Cells.Item(3, 19).Value = 0.2372736994476142
$ws.Cells.Item(3, 20).Value = 0.2372736994476143
$ws.Cells.Item(4, 7).Value = 22.77954133333333
$ws.Cells.Item(4, 8).Value = 68.338624
$ws.Cells.Item(4, 9).Value = 0.8649343844704168
$ws.Cells.Item(4, 10).Value = 0.8649343844704167
$ws.Cells.Item(4, 13).Value = 39.361408
$ws.Cells.Item(4, 14).Value = 118.084224
$ws.Cells.Item(4, 15).Value = 0.3925813724534833
$ws.Cells.Item(4, 16).Value = 0.3925813724534833
$ws.Cells.Item(4, 17).Value = 896.6348204741973
$ws.Cells.Item(4, 18).Value = 8069.713384267776
$ws.Cells.Item(4, 19).Value = 0.339557127737605
$ws.Cells.Item(4, 20).Value = 0.339557127737605
$ws.Cells.Item(5, 7).Value = 22.77954133333333
$ws.Cells.Item(5, 8).Value = 68.338624
$ws.Cells.Item(5, 9).Value = 0.8649343844704168
$ws.Cells.Item(5, 10).Value = 0.8649343844704167
$ws.Cells.Item(5, 13).Value = 32.221985
$ws.Cells.Item(5, 14).Value = 96.665955
$ws.Cells.Item(5, 15).Value = 0.3213744562815322
$ws.Cells.Item(5, 16).Value = 0.3213744562815322
$ws.Cells.Item(5, 17).Value = 734.0020391495465
$ws.Cells.Item(5, 18).Value = 6606.018352345919
$ws.Cells.Item(5, 19).Value = 0.2779678175283819
$ws.Cells.Item(5, 20).Value = 0.2779678175283819
$ws.Cells.Item(6, 9).Value = 0.008798055815159926
$ws.Cells.Item(6, 10).Value = 0.008798055815159925
$ws.Cells.Item(6, 13).Value = 1.174933333333333
$ws.Cells.Item(6, 14).Value = 3.5248
$ws.Cells.Item(6, 15).Value = 0.01171850713626266
$ws.Cells.Item(6, 16).Value = 0.01171850713626266
$ws.Cells.Item(6, 17).Value = 0.2722461525333333
$ws.Cells.Item(6, 18).Value = 2.4502153728
$ws.Cells.Item(6, 19).Value = 0.0001031000798551888
$ws.Cells.Item(6, 20).Value = 0.0001031000798551887
$ws.Cells.Item(7, 9).Value = 0.008798055815159926
$ws.Cells.Item(7, 10).Value = 0.008798055815159925
$ws.Cells.Item(7, 15).Value = 0.2743256641287217
$ws.Cells.Item(7, 16).Value = 0.2743256641287218
$ws.Cells.Item(7, 19).Value = 0.002413532504535309
$ws.Cells.Item(7, 20).Value = 0.002413532504535309
$ws.Cells.Item(8, 9).Value = 0.008798055815159926
$ws.Cells.Item(8, 10).Value = 0.008798055815159925
$ws.Cells.Item(8, 13).Value = 39.361408
$ws.Cells.Item(8, 14).Value = 118.084224
$ws.Cells.Item(8, 15).Value = 0.3925813724534833
$ws.Cells.Item(8, 16).Value = 0.3925813724534833
$ws.Cells.Item(8, 17).Value = 9.120510570496002
$ws.Cells.Item(8, 18).Value = 82.084595134464
$ws.Cells.Item(8, 19).Value = 0.003453952826837834
$ws.Cells.Item(8, 20).Value = 0.003453952826837834
$ws.Cells.Item(9, 9).Value = 0.008798055815159926
$ws.Cells.Item(9, 10).Value = 0.008798055815159925
$ws.Cells.Item(9, 13).Value = 32.221985
$ws.Cells.Item(9, 14).Value = 96.665955
$ws.Cells.Item(9, 15).Value = 0.3213744562815322
$ws.Cells.Item(9, 16).Value = 0.3213744562815322
$ws.Cells.Item(9, 17).Value = 7.46622058832
$ws.Cells.Item(9, 18).Value = 67.19598529487999
$ws.Cells.Item(9, 19).Value = 0.002827470403931594
$ws.Cells.Item(9, 20).Value = 0.002827470403931593
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.07785033333333334
$ws.Cells.Item(10, 8).Value = 0.233551
$ws.Cells.Item(10, 9).Value = 0.002955960752552617
$ws.Cells.Item(10, 10).Value = 0.002955960752552617
$ws.Cells.Item(10, 13).Value = 1.174933333333333
$ws.Cells.Item(10, 14).Value = 3.5248
$ws.Cells.Item(10, 15).Value = 0.01171850713626266
$ws.Cells.Item(10, 16).Value = 0.01171850713626266
$ws.Cells.Item(10, 17).Value = 0.09146895164444446
$ws.Cells.Item(10, 18).Value = 0.8232205648
$ws.Cells.Item(10, 19).Value = 0.00003463944717330018
$ws.Cells.Item(10, 20).Value = 0.00003463944717330017
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.07785033333333334
$ws.Cells.Item(11, 8).Value = 0.233551
$ws.Cells.Item(11, 9).Value = 0.002955960752552617
$ws.Cells.Item(11, 10).Value = 0.002955960752552617
$ws.Cells.Item(11, 15).Value = 0.2743256641287217
$ws.Cells.Item(11, 16).Value = 0.2743256641287218
$ws.Cells.Item(11, 17).Value = 2.141252346843111
$ws.Cells.Item(11, 18).Value = 19.271271121588
$ws.Cells.Item(11, 19).Value = 0.0008108958965824327
$ws.Cells.Item(11, 20).Value = 0.0008108958965824328
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.07785033333333334
$ws.Cells.Item(12, 8).Value = 0.233551
$ws.Cells.Item(12, 9).Value = 0.002955960752552617
$ws.Cells.Item(12, 10).Value = 0.002955960752552617
$ws.Cells.Item(12, 13).Value = 39.361408
$ws.Cells.Item(12, 14).Value = 118.084224
$ws.Cells.Item(12, 15).Value = 0.3925813724534833
$ws.Cells.Item(12, 16).Value = 0.3925813724534833
$ws.Cells.Item(12, 17).Value = 3.064298733269334
$ws.Cells.Item(12, 18).Value = 27.578688599424
$ws.Cells.Item(12, 19).Value = 0.001160455129155738
$ws.Cells.Item(12, 20).Value = 0.001160455129155738
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.07785033333333334
$ws.Cells.Item(13, 8).Value = 0.233551
$ws.Cells.Item(13, 9).Value = 0.002955960752552617
$ws.Cells.Item(13, 10).Value = 0.002955960752552617
$ws.Cells.Item(13, 13).Value = 32.221985
$ws.Cells.Item(13, 14).Value = 96.665955
$ws.Cells.Item(13, 15).Value = 0.3213744562815322
$ws.Cells.Item(13, 16).Value = 0.3213744562815322
$ws.Cells.Item(13, 17).Value = 2.508492272911667
$ws.Cells.Item(13, 18).Value = 22.576430456205
$ws.Cells.Item(13, 19).Value = 0.0009499702796411459
$ws.Cells.Item(13, 20).Value = 0.000949970279641146
$ws.Cells.Item(14, 7).Value = 3.247624
$ws.Cells.Item(14, 8).Value = 9.742872
$ws.Cells.Item(14, 9).Value = 0.1233115989618705
$ws.Cells.Item(14, 10).Value = 0.1233115989618705
$ws.Cells.Item(14, 13).Value = 1.174933333333333
$ws.Cells.Item(14, 14).Value = 3.5248
$ws.Cells.Item(14, 15).Value = 0.01171850713626266
$ws.Cells.Item(14, 16).Value = 0.01171850713626266
$ws.Cells.Item(14, 17).Value = 3.815741691733334
$ws.Cells.Item(14, 18).Value = 34.3416752256
$ws.Cells.Item(14, 19).Value = 0.001445027852418638
$ws.Cells.Item(14, 20).Value = 0.001445027852418638
$ws.Cells.Item(15, 7).Value = 3.247624
$ws.Cells.Item(15, 8).Value = 9.742872
$ws.Cells.Item(15, 9).Value = 0.1233115989618705
$ws.Cells.Item(15, 10).Value = 0.1233115989618705
$ws.Cells.Item(15, 15).Value = 0.2743256641287217
$ws.Cells.Item(15, 16).Value = 0.2743256641287218
$ws.Cells.Item(15, 17).Value = 89.32501909643733
$ws.Cells.Item(15, 18).Value = 803.9251718679361
$ws.Cells.Item(15, 19).Value = 0.03382753627998972
$ws.Cells.Item(15, 20).Value = 0.03382753627998972
$ws.Cells.Item(16, 7).Value = 3.247624
$ws.Cells.Item(16, 8).Value = 9.742872
$ws.Cells.Item(16, 9).Value = 0.1233115989618705
$ws.Cells.Item(16, 10).Value = 0.1233115989618705
$ws.Cells.Item(16, 13).Value = 39.361408
$ws.Cells.Item(16, 14).Value = 118.084224
$ws.Cells.Item(16, 15).Value = 0.3925813724534833
$ws.Cells.Item(16, 16).Value = 0.3925813724534833
$ws.Cells.Item(16, 17).Value = 127.831053294592
$ws.Cells.Item(16, 18).Value = 1150.479479651328
$ws.Cells.Item(16, 19).Value = 0.04840983675988466
$ws.Cells.Item(16, 20).Value = 0.04840983675988465
$ws.Cells.Item(17, 7).Value = 3.247624
$ws.Cells.Item(17, 8).Value = 9.742872
$ws.Cells.Item(17, 9).Value = 0.1233115989618705
$ws.Cells.Item(17, 10).Value = 0.1233115989618705
$ws.Cells.Item(17, 13).Value = 32.221985
$ws.Cells.Item(17, 14).Value = 96.665955
$ws.Cells.Item(17, 15).Value = 0.3213744562815322
$ws.Cells.Item(17, 16).Value = 0.3213744562815322
$ws.Cells.Item(17, 17).Value = 104.64489181364
$ws.Cells.Item(17, 18).Value = 941.80402632276
$ws.Cells.Item(17, 19).Value = 0.03962919806957749
$ws.Cells.Item(17, 20).Value = 0.03962919806957749
